# Update ObjTables metadata strings embedded in row/column A1 (and A2 of
# the table-of-contents sheet):
#   - bump objTablesVersion 0.0.8 -> 0.0.9
#   - bump date 2020-03-14 11:51:03 -> 2020-04-27 01:05:04
#   - rename the Data-table header attribute id= to class=
#
# The target cells live on sheets that have sheet-level protection enabled
# (sheetProtection sheet="1" objects="1" scenarios="1" insertRows="0"
# deleteRows="0"). Calling Worksheet.Protect()/Unprotect() in this runtime
# only round-trips the basic "sheet" flag and drops the other attributes,
# so instead we temporarily flip Locked off on just the cells we need to
# touch (which Excel allows on a protected sheet), make the edits, and
# flip Locked back on - this leaves the sheet's <sheetProtection> element
# completely untouched.

$wb = $excel.ActiveWorkbook

$tocSheet    = $wb.Worksheets.Item("!!_Table of contents")
$schemaSheet = $wb.Worksheets.Item("!!_Schema")
$childSheet  = $wb.Worksheets.Item("!!Child")
$parentSheet = $wb.Worksheets.Item("!!Parent")

$targets = @(
    @{ ws = $tocSheet;    addr = "A1" },
    @{ ws = $tocSheet;    addr = "A2" },
    @{ ws = $schemaSheet; addr = "A1" },
    @{ ws = $childSheet;  addr = "A1" },
    @{ ws = $parentSheet; addr = "A1" }
)

foreach ($t in $targets) {
    $t.ws.Range($t.addr).Locked = $false
}

$tocSheet.Range("A1").Value    = "!!!ObjTables objTablesVersion='0.0.9' date='2020-04-27 01:05:04'"
$tocSheet.Range("A2").Value    = "!!ObjTables type='TableOfContents' tableFormat='row' description='Table of contents' date='2020-04-27 01:05:04' objTablesVersion='0.0.9'"
$schemaSheet.Range("A1").Value = "!!ObjTables type='Schema' tableFormat='row' description='Table/model and column/attribute definitions' date='2020-04-27 01:05:04' objTablesVersion='0.0.9'"
$childSheet.Range("A1").Value  = "!!ObjTables type='Data' tableFormat='row' class='Child' name='Child' date='2020-04-27 01:05:04' objTablesVersion='0.0.9'"
$parentSheet.Range("A1").Value = "!!ObjTables type='Data' tableFormat='column' class='Parent' name='Parent' date='2020-04-27 01:05:04' objTablesVersion='0.0.9'"

foreach ($t in $targets) {
    $t.ws.Range($t.addr).Locked = $true
}
